$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix two cells that had a stray hard line-break inside the text (the
#    wrapped text is joined back into a single line with a plain space).
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 3).Value = "VGR des Bundes - Bruttowertschöpfung, Bruttoinlandsprodukt (nominal/preisbereinigt)"
$ws.Cells.Item(23, 4).Value = "National accounts - Gross value added, gross domestic product"

# ---------------------------------------------------------------------------
# 2) Update the Gini-coefficient row (row 39): shorten the German/English
#    descriptions and refresh the Eurostat links.
# ---------------------------------------------------------------------------
$ws.Cells.Item(39, 3).Value = "Gini-Koeffizient des verfügbaren Äquivalenzeinkommens vor Sozialleistungen - EU-SILC Erhebung"
$ws.Cells.Item(39, 4).Value = "Gini coefficient of equivalised disposable income before social transfers"
$ws.Cells.Item(39, 5).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_DI12C/default/table?lang=de&category=livcon.ilc.ilc_ie.ilc_iei"
$ws.Cells.Item(39, 6).Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_DI12C/default/table?lang=en&category=livcon.ilc.ilc_ie.ilc_iei"

# ---------------------------------------------------------------------------
# 3) Append seven new source rows (69-75), matching the look & feel
#    (style/formatting) of the existing data rows by copying the format of
#    the last existing row before filling in the new values.
# ---------------------------------------------------------------------------
$ws.Range("A68:F68").Copy()
$ws.Range("A69:F75").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    @("L_RDB_1", "Q_DESTATIS", "Bevölkerung", "Population (only available in German)", "https://www.regionalstatistik.de/genesis/online?operation=previous&levelindex=0&step=0&titel=Tabellenaufbau&levelid=1668672879939&acceptscookies=false#abreadcrumb", ""),
    @("L_RDB_2", "Q_DESTATIS", "Siedlungsfläche", "Settlement area (only available in German)", "https://www.regionalstatistik.de/genesis//online?operation=table&code=33111-02-01-4&bypass=true&levelindex=1&levelid=1668501420853#abreadcrumb", ""),
    @("L_RDB_3", "Q_DESTATIS", "Verkehrsfläche", "Transport area (only available in German)", "https://www.regionalstatistik.de/genesis//online?operation=table&code=33111-03-01-4&bypass=true&levelindex=1&levelid=1668501420853#abreadcrumb", ""),
    @("L_IFEU_1", "Q_IFEU", "TREMOD", "TREMOD", "https://www.ifeu.de/methoden-tools/modelle/tremod/", "https://www.ifeu.de/en/methods-tools/models/tremod/"),
    @("L_ERSTT_11", "Q_EUROSTAT", "Quote der Überbelastung durch Wohnkosten", "Housing cost overburden rate", "https://ec.europa.eu/eurostat/databrowser/view/ILC_LVHO07A/default/table?lang=de&category=livcon.ilc.ilc_lv.ilc_lvho.ilc_lvho_hc", "https://ec.europa.eu/eurostat/databrowser/view/ILC_LVHO07A/default/table?category=livcon.ilc.ilc_lv.ilc_lvho.ilc_lvho_hc"),
    @("L_DSTTS_24", "Q_DESTATIS", "Überbelastung durch Wohnkosten", "Housing cost overburden", "https://www.destatis.de/Europa/DE/Thema/Bevoelkerung-Arbeit-Soziales/Soziales-Lebensbedingungen/Wohnkosten.html", "https://www.destatis.de/Europa/EN/Topic/Population-Labour-Social-Issues/Social-issues-living-conditions/_node.html;jsessionid=B340DD00C6EEDC7477B2AD2B54E4BC40.live731#587120"),
    @("L_DSTTS_25", "Q_DESTATIS", "Umweltökonomische Gesamtrechnungen", "Environmental Economic Accounting", "https://www.destatis.de/DE/Themen/Gesellschaft-Umwelt/Umwelt/UGR/_inhalt.html", "https://www.destatis.de/EN/Themes/Society-Environment/Environmental-Economic-Accounting/_node.html")
)

$startRow = 69
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $data[$col - 1]
    }
}
